# Edit script: insert 3 new price rows for "Durazno" at Vega Monumental
# Concepción (Carson / Elegant Lady / Polar King, "Primera" quality,
# volume 220, date 2022-01-28) right above the existing "Kakamas" block
# that currently starts at row 82. This pushes the existing rows 82:152
# down to 85:155 and grows the sheet from 152 to 155 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at the top of the block (row 82), shifting the
# existing data (old rows 82-152) down to rows 85-155.
$ws.Rows("82:84").Insert()

# Values that stay constant across this whole data block.
$marketId   = 11
$market     = "Vega Monumental Concepción"
$region     = "Bíobío"
$codreg     = 8
$tipo       = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103004
$categoria  = "Durazno"
$unidad     = "$/caja 16 kilos empedrada"
$origen     = "Región de O'Higgins"
$fecha      = 44589

$newRows = @(
    @{ Row = 82; Variedad = "Carson";       Calidad = "Primera"; Volumen = 220; PMin = 13000; PMax = 14000; PProm = 13545; PKg = 847; KgUnidad = 16 },
    @{ Row = 83; Variedad = "Elegant Lady";  Calidad = "Primera"; Volumen = 220; PMin = 12000; PMax = 13000; PProm = 12455; PKg = 778; KgUnidad = 16 },
    @{ Row = 84; Variedad = "Polar King";    Calidad = "Primera"; Volumen = 220; PMin = 12000; PMax = 13000; PProm = 12455; PKg = 778; KgUnidad = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $marketId
    $ws.Cells.Item($row, 2).Value  = $market
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
